$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 / Row 13 value updates (labels/structure unchanged, text content changes) ---
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("B13").Value = "6279110 - Carlos Alberto Moreira dos Santos"
$ws.Range("C13").Value = "6279110 - Carlos Alberto Moreira dos Santos"

# --- Remove old rows 14-30 entirely (they will be rebuilt from scratch) ---
$ws.Range("A14:A30").EntireRow.Delete()

# --- Rebuild rows 14-26 with the new content, pulling cell formats from
#     known-good template cells so styles (s=1/2/3) match exactly ---
$ws.Range("A10:C10").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C14").Value = "6495737 - Durval Rodrigues Junior"
$ws.Rows(14).RowHeight = 60

$ws.Range("A12").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Short syllabus:"
$ws.Rows(15).RowHeight = 60

$ws.Range("A10:C10").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C16").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows(16).RowHeight = 120

$ws.Range("A12").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows(17).RowHeight = 120

$ws.Range("A12").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A10:C10").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C19").Value = "1341653 - Maria José Ramos Sandim"
$ws.Rows(19).RowHeight = 60

$ws.Range("A10:C10").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C20").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Rows(20).RowHeight = 60

$ws.Range("A10:C10").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C21").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Rows(21).RowHeight = 60

$ws.Range("A10:C10").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows(22).RowHeight = 120

$ws.Range("A12").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B13:C13").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Range("B24").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("C24").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Rows(24).RowHeight = 30

$ws.Range("B13:C13").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$ws.Range("B25").Value = "LOM3253 -  Física Matemática  (Requisito)`n"
$ws.Range("C25").Value = "LOM3253 -  Física Matemática  (Requisito)`n"
$ws.Rows(25).RowHeight = 30

$ws.Range("B13:C13").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("B26").Value = "LOM3257 -  Mecânica Clássica  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOM3257 -  Mecânica Clássica  (Requisito fraco)`n"
$ws.Rows(26).RowHeight = 30

$excel.CutCopyMode = 0
$ws.Range("A1").Select()
